# Weekly update: insert a new "Camote/Paine" record as the latest week
# for the Zapallo (Hortaliza) sheet, pushing every existing record down
# by one row. This mirrors a new row being prepended to the dataset.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 73..106 down to 74..107, opening up a blank row 73.
$ws.Rows.Item(73).Insert()

# Populate the newly-opened row 73 with the new weekly record.
$ws.Cells.Item(73, 1).Value  = 7
$ws.Cells.Item(73, 2).Value  = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(73, 3).Value  = "Ñuble"
$ws.Cells.Item(73, 4).Value  = 44523
$ws.Cells.Item(73, 5).Value  = 16
$ws.Cells.Item(73, 6).Value  = 100112045
$ws.Cells.Item(73, 7).Value  = "Zapallo"
$ws.Cells.Item(73, 8).Value  = "Paine"
$ws.Cells.Item(73, 9).Value  = "1a (guarda)"
$ws.Cells.Item(73, 10).Value = 200
$ws.Cells.Item(73, 11).Value = 220
$ws.Cells.Item(73, 12).Value = 250
$ws.Cells.Item(73, 13).Value = 235
$ws.Cells.Item(73, 14).Value = "$/kilo (volumen en unidades)"
$ws.Cells.Item(73, 15).Value = "Región de O'Higgins"
$ws.Cells.Item(73, 16).Value = 235
$ws.Cells.Item(73, 17).Value = 1
$ws.Cells.Item(73, 18).Value = "Hortaliza"
